$d = $word.ActiveDocument

# The document currently looks like:
#   1. Title paragraph: "Write Up"
#   2. Empty paragraph (Normal style)
#   3. Empty paragraph (Heading 1 style)
#
# We need to insert three new paragraphs right before the existing empty
# (Normal-style) paragraph #2, leaving paragraphs #2 and #3 untouched:
#   - "This week, we will be taking look at creating the room. ..."
#   - "So, if this sounds like something that you would like to learn ..."
#   - "4 The Room" (Heading 1 style)

# Locate the empty paragraph that currently follows the Title paragraph.
$anchor = $d.Paragraphs(2)
$r = $anchor.Range
$r.Collapse(1)   # wdCollapseStart

# Reserve three new empty paragraphs before the anchor paragraph. Using
# InsertParagraphBefore keeps the anchor paragraph's (Normal, no explicit
# style) formatting for the newly created paragraphs.
$r.InsertParagraphBefore()
$r.InsertParagraphBefore()
$r.InsertParagraphBefore()

# Fill in the text for the three new paragraphs.
$d.Paragraphs(2).Range.Text = "This week, we will be taking look at creating the room. A room is needed in order to have a game. So, it is and extremely important concept in Game Maker. Learn to create layers and then also how to bring your game elements into that room."
$d.Paragraphs(3).Range.Text = "So, if this sounds like something that you would like to learn a bit more about then please join us for our brand-new article entitled:"
$d.Paragraphs(4).Range.Text = "4 The Room"
$d.Paragraphs(4).Style = "Heading 1"
